$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the original
# "t=inlineStr"-like plain-text representation instead of letting Excel's
# auto-detection coerce numeric-looking strings (e.g. "1.000", "5.860")
# into real numbers. We briefly force Text number format for the write,
# then restore General/Normal so no stray style gets attached to the cell.
function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Rows 13 and 14 swap places: the coin that used to be row 13 (WrappedEther)
# is now row 14, and the coin that used to be row 14 (Polkadot) is now row
# 13 - each with refreshed price/volume figures.
Set-TextValue $ws "B13" "Polkadot"
Set-TextValue $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D13" "5.860"
Set-TextValue $ws "E13" "  -4.15%  "

Set-TextValue $ws "B14" "WrappedEther"
Set-TextValue $ws "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D14" "1.821.97"
Set-TextValue $ws "E14" "  -4.28%  "

# Refreshed Price / Volume(1h) figures for the remaining rows.
Set-TextValue $ws "D2" "28.656.29"
Set-TextValue $ws "E2" "  -3.10%  "
Set-TextValue $ws "D3" "1.853.12"
Set-TextValue $ws "E3" "  -3.64%  "
Set-TextValue $ws "D5" "334.82"
Set-TextValue $ws "E5" "  +2.58%  "
Set-TextValue $ws "D6" "1.000"
Set-TextValue $ws "E6" "  -1.11%  "
Set-TextValue $ws "D7" "0.4657"
Set-TextValue $ws "E7" "  -3.25%  "
Set-TextValue $ws "D8" "0.3914"
Set-TextValue $ws "E8" "  -3.45%  "
Set-TextValue $ws "D9" "46.15"
Set-TextValue $ws "E9" "  -3.37%  "
Set-TextValue $ws "D10" "0.07906"
Set-TextValue $ws "E10" "  -3.79%  "
Set-TextValue $ws "D11" "0.9851"
Set-TextValue $ws "E11" "  -2.51%  "
Set-TextValue $ws "D12" "22.35"
Set-TextValue $ws "E12" "  -6.35%  "
Set-TextValue $ws "D15" "7.015"
Set-TextValue $ws "E15" "  -3.84%  "
Set-TextValue $ws "D16" "0.06877"
Set-TextValue $ws "E16" "  +0.13%  "
Set-TextValue $ws "D17" "87.72"
Set-TextValue $ws "E17" "  -4.40%  "
Set-TextValue $ws "D18" "1.001"
Set-TextValue $ws "E18" "  -1.21%  "
Set-TextValue $ws "E19" "  -3.21%  "
Set-TextValue $ws "D20" "17.17"
Set-TextValue $ws "E20" "  -2.64%  "
Set-TextValue $ws "E21" "  -1.01%  "
Set-TextValue $ws "D22" "28.670.37"
Set-TextValue $ws "E22" "  -3.07%  "
Set-TextValue $ws "D23" "5.406"
Set-TextValue $ws "E23" "  -4.74%  "
Set-TextValue $ws "D24" "11.36"
Set-TextValue $ws "E24" "  -5.37%  "
Set-TextValue $ws "D25" "2.142"
Set-TextValue $ws "E25" "  -2.15%  "
Set-TextValue $ws "D26" "2.021.81"
Set-TextValue $ws "E26" "  -5.74%  "
Set-TextValue $ws "D27" "153.43"
Set-TextValue $ws "E27" "  -1.80%  "
Set-TextValue $ws "D28" "19.52"
Set-TextValue $ws "E28" "  -2.65%  "
Set-TextValue $ws "D29" "6.055"
Set-TextValue $ws "E29" "  -5.46%  "
Set-TextValue $ws "D30" "2.033"
Set-TextValue $ws "E30" "  -2.76%  "
Set-TextValue $ws "D31" "117.76"
Set-TextValue $ws "E31" "  -2.47%  "
Set-TextValue $ws "D32" "0.9781"
Set-TextValue $ws "E32" "  -3.42%  "
Set-TextValue $ws "D33" "0.09397"
Set-TextValue $ws "E33" "  -2.14%  "
Set-TextValue $ws "D35" "3.482"
Set-TextValue $ws "E35" "  -2.41%  "
Set-TextValue $ws "D36" "1.351"
Set-TextValue $ws "E36" "  -2.02%  "
Set-TextValue $ws "D37" "0.06161"
Set-TextValue $ws "E37" "  -3.09%  "
Set-TextValue $ws "D38" "0.02201"
Set-TextValue $ws "E38" "  -3.79%  "
Set-TextValue $ws "D39" "1.166"
Set-TextValue $ws "E39" "  -2.09%  "
Set-TextValue $ws "D40" "0.5734"
Set-TextValue $ws "E40" "  -3.67%  "
Set-TextValue $ws "D41" "7.637"
Set-TextValue $ws "E41" "  -2.92%  "
Set-TextValue $ws "D42" "10.21"
Set-TextValue $ws "E42" "  -5.00%  "
Set-TextValue $ws "D43" "0.1803"
Set-TextValue $ws "E43" "  -2.59%  "
Set-TextValue $ws "D44" "2.367"
Set-TextValue $ws "E44" "  -2.53%  "
Set-TextValue $ws "D45" "1.229"
Set-TextValue $ws "E45" "  -1.65%  "
Set-TextValue $ws "D46" "0.5402"
Set-TextValue $ws "E46" "  -2.79%  "
Set-TextValue $ws "D47" "11.75"
Set-TextValue $ws "E47" "  -5.58%  "
Set-TextValue $ws "D48" "0.07155"
Set-TextValue $ws "E48" "  -5.23%  "
Set-TextValue $ws "D49" "1.915"
Set-TextValue $ws "E49" "  -3.89%  "
Set-TextValue $ws "D50" "115.66"
Set-TextValue $ws "E50" "  -3.22%  "
Set-TextValue $ws "D51" "43.41"
Set-TextValue $ws "E51" "  +2.46%  "
